$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9370161890983582
$ws.Range("B1").Value = 1.5932697057724
$ws.Range("C1").Value = 3.19239068031311
$ws.Range("D1").Value = 3.08130407333374
$ws.Range("E1").Value = 0.3484648764133453
